# Cluster: import wizard with video list
#
# 1. Rename the sheet from "Sheet1" to "cluster".
# 2. Add a new "video_codes" column (D) with a header styled like the
#    existing header row, and populate it for the clusters that have
#    associated videos (clusters 07, 12 and 13 have none).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "cluster"

# --- Header -----------------------------------------------------------
$ws.Range("D1").Value = "video_codes"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null

# --- Data rows ----------------------------------------------------------
$videoCodes = @{
    2  = "vdo_0033;vdo_0035;vdo_0038"
    3  = "vdo_0023;vdo_0031;vdo_0043"
    4  = "vdo_0007;vdo_0011;vdo_0012;vdo_0015;vdo_0019;vdo_0045;vdo_0047;vdo_0049;vdo_0050;vdo_0051"
    5  = "vdo_0013;vdo_0017;vdo_0026;vdo_0046"
    6  = "vdo_0027;vdo_0048;vdo_0052;vdo_0053;vdo_0054;vdo_0055;vdo_0056;vdo_0057;vdo_0058;vdo_0059"
    7  = "vdo_0021;vdo_0030"
    9  = "vdo_0034;vdo_0036;vdo_0037;vdo_0041"
    10 = "vdo_0004;vdo_0024;vdo_0025;vdo_0028"
    11 = "vdo_0009;vdo_0010;vdo_0029;vdo_0039;vdo_0044"
    12 = "vdo_0002;vdo_0005;vdo_0006;vdo_0008;vdo_0014;vdo_0022;vdo_0040;vdo_0042"
    15 = "vdo_0001;vdo_0003"
    16 = "vdo_0016;vdo_0018;vdo_0020"
    17 = "vdo_0032"
}

$ws.Range("A2").Copy() | Out-Null

foreach ($row in $videoCodes.Keys | Sort-Object) {
    $cell = $ws.Range("D$row")
    $cell.Value = $videoCodes[$row]
    $ws.Range("D$row").PasteSpecial(-4122) | Out-Null
}
